# Name-change order template fix:
# Wrap the p.birth_state placeholder with a state_name() call so the
# rendered order prints the full state name instead of the raw
# abbreviation/code, e.g.:
#   ... in the state of {{ p.birth_state }} ...
# becomes
#   ... in the state of {{ state_name(p.birth_state) }} ...

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "p.birth_state",       # FindText
    $false,                 # MatchCase
    $false,                 # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                 # Format
    "state_name(p.birth_state)",  # ReplaceWith
    2                        # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find 'p.birth_state' placeholder to update."
}
